$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2: I2:L2
$ws.Range("I2").Value = 0.213557390838831
$ws.Range("J2").Value = 0.2800502618304645
$ws.Range("K2").Value = -0.291346792213953
$ws.Range("L2").Value = 2.49331058952106

# Row 20: I20:L20
$ws.Range("I20").Value = -0.003370232952257873
$ws.Range("J20").Value = 0.381881807301952
$ws.Range("K20").Value = 0.01308274625048825
$ws.Range("L20").Value = 1.940217643963906
